$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => @(NewPrice-or-$null, NewVolume)
$updates = @{
    2  = @("38.749.85", "  +1.45%  ")
    3  = @("2.099.39",  "  +0.31%  ")
    4  = @($null,       "  -0.05%  ")
    5  = @("228.95",    "  +0.16%  ")
    6  = @($null,       "  +0.85%  ")
    7  = @("61.51",     "  +0.79%  ")
    8  = @($null,       "  -0.06%  ")
    9  = @($null,       "  +1.89%  ")
    10 = @($null,       "  -0.67%  ")
    11 = @($null,       "  +0.13%  ")
    12 = @("15.57",     "  +6.01%  ")
    13 = @("2.409.82",  "  +0.30%  ")
    14 = @("21.94",     "  -1.47%  ")
    15 = @($null,       "  +3.92%  ")
    16 = @($null,       "  +0.34%  ")
    17 = @("2.095.99",  "  -0.51%  ")
    18 = @("38.766.85", "  +1.72%  ")
    19 = @($null,       "  +2.54%  ")
    20 = @("6.07",      "  +0.74%  ")
    21 = @($null,       "  +0.36%  ")
    22 = @("227.55",    "  +1.55%  ")
    24 = @($null,       "  -2.24%  ")
    26 = @("172.03",    "  +1.19%  ")
    27 = @("9.56",      "  +1.02%  ")
    28 = @($null,       "  +5.51%  ")
    29 = @("1.41",      "  +3.81%  ")
    30 = @("19.33",     "  +1.92%  ")
    31 = @("2.48",      "  +3.73%  ")
    32 = @($null,       "  +0.97%  ")
    33 = @($null,       "  +2.19%  ")
    34 = @("4.77",      "  +1.54%  ")
    36 = @("6.59",      "  +2.92%  ")
    37 = @("2.41",      "  +0.79%  ")
    38 = @("3.56",      "  +1.23%  ")
    39 = @("0.999",     "  -0.04%  ")
    40 = @("18.26",     "  +0.88%  ")
    41 = @($null,       "  +4.41%  ")
    42 = @("101.53",    "  +1.42%  ")
    43 = @("1.533.10",  "  -1.35%  ")
    44 = @($null,       "  -0.96%  ")
    45 = @($null,       "  +4.04%  ")
    46 = @("0.0911",    "  -0.22%  ")
    47 = @("1.13",      "  +1.58%  ")
    48 = @($null,       "  -1.00%  ")
    49 = @($null,       "  +1.80%  ")
    50 = @($null,       "  -0.93%  ")
    51 = @("2.292.60",  "  +0.15%  ")
}

# Rows whose new Price string parses as a plain number (no thousands-dot
# grouping) and therefore needs the cell pre-formatted as Text - otherwise
# Excel's type auto-detection would silently store it as a Double instead of
# keeping it a string like the rest of the "Price" column.
$needsTextForce = @(5, 7, 12, 14, 20, 22, 26, 27, 29, 30, 31, 34, 36, 37, 38, 39, 40, 42, 46, 47)

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $newPrice = $pair[0]
    $newVolume = $pair[1]

    if ($null -ne $newPrice) {
        $priceCell = $ws.Range("D$row")
        if ($needsTextForce -contains $row) {
            $priceCell.NumberFormat = "@"
        }
        $priceCell.Value = $newPrice
    }
    $ws.Range("E$row").Value = $newVolume
}
